$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")

$rngA2 = $wsAbout.Range("A2")
$rngA2.Value = $rngA2.Value().Replace($oldStamp, $newStamp)

$rngA6 = $wsAbout.Range("A6")
$rngA6.Value = $rngA6.Value().Replace($oldStamp, $newStamp)

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 9; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19 -> build_version
    $cell.Value = $cell.Value().Replace($oldStamp, $newStamp)
}
